$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.975.88'
$ws.Range('E2').Value = '  +6.51%  '
$ws.Range('D3').Value = '3.669.16'
$ws.Range('E3').Value = '  +18.96%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.57'
$ws.Range('E5').Value = '  +4.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.45'
$ws.Range('E6').Value = '  +7.36%  '
$ws.Range('D7').Value = '3.666.85'
$ws.Range('E7').Value = '  +18.98%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.535'
$ws.Range('E9').Value = '  +4.55%  '
$ws.Range('E10').Value = '  +7.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.53'
$ws.Range('E11').Value = '  +3.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.501'
$ws.Range('E12').Value = '  +7.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.99'
$ws.Range('E13').Value = '  +11.52%  '
$ws.Range('E14').Value = '  +6.25%  '
$ws.Range('D15').Value = '4.286.12'
$ws.Range('E15').Value = '  +19.08%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '70.991.73'
$ws.Range('E16').Value = '  +6.63%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.672.08'
$ws.Range('E17').Value = '  +19.08%  '
$ws.Range('E18').Value = '  +2.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.50'
$ws.Range('E19').Value = '  +8.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.01'
$ws.Range('E20').Value = '  +2.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '519.15'
$ws.Range('E21').Value = '  +7.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.16'
$ws.Range('E22').Value = '  +16.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.744'
$ws.Range('E23').Value = '  +8.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.91'
$ws.Range('E24').Value = '  +5.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.46'
$ws.Range('E25').Value = '  +6.59%  '
$ws.Range('E26').Value = '  +8.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.85'
$ws.Range('E27').Value = '  +8.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.53'
$ws.Range('E29').Value = '  +12.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.15'
$ws.Range('E30').Value = '  +3.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.68'
$ws.Range('E31').Value = '  +13.43%  '
$ws.Range('E32').Value = '  +6.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0000109'
$ws.Range('E33').Value = '  +17.90%  '
$ws.Range('E34').Value = '  +4.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.19'
$ws.Range('E36').Value = '  +10.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.02'
$ws.Range('E37').Value = '  +8.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.342'
$ws.Range('E38').Value = '  +12.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.13'
$ws.Range('E39').Value = '  +8.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.70'
$ws.Range('E40').Value = '  +3.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '46.37'
$ws.Range('E41').Value = '  -5.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.127'
$ws.Range('E42').Value = '  +4.20%  '
$ws.Range('D43').Value = '3.177.70'
$ws.Range('E43').Value = '  +14.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.85'
$ws.Range('E44').Value = '  +7.61%  '
$ws.Range('E45').Value = '  +8.27%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '398.73'
$ws.Range('E46').Value = '  +9.18%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0368'
$ws.Range('E47').Value = '  +7.55%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '28.24'
$ws.Range('E48').Value = '  +16.20%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '135.69'
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.45'
$ws.Range('E51').Value = '  +12.78%  '
